$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8..22 hold the three region blocks (Brasil/Nordeste/Sergipe) x years.
# A new 2024 data point was inserted for each region, shifting the
# Nordeste and Sergipe blocks down by one row and adding two new rows
# at the bottom for Sergipe.

# Make sure column C (dates-as-text) is never auto-converted to a real
# date serial by Excel's type inference - force Text format first.
$ws.Range("C8:C22").NumberFormat = "@"

$data = @(
    @("Brasil",   "Número médio de moradores", "01/01/2024", 2.7),
    @("Nordeste", "Número médio de moradores", "01/01/2016", 3.2),
    @("Nordeste", "Número médio de moradores", "01/01/2017", 3.2),
    @("Nordeste", "Número médio de moradores", "01/01/2018", 3.1),
    @("Nordeste", "Número médio de moradores", "01/01/2019", 3.1),
    @("Nordeste", "Número médio de moradores", "01/01/2022", 3),
    @("Nordeste", "Número médio de moradores", "01/01/2023", 2.8),
    @("Nordeste", "Número médio de moradores", "01/01/2024", 2.8),
    @("Sergipe",  "Número médio de moradores", "01/01/2016", 3),
    @("Sergipe",  "Número médio de moradores", "01/01/2017", 3),
    @("Sergipe",  "Número médio de moradores", "01/01/2018", 3.1),
    @("Sergipe",  "Número médio de moradores", "01/01/2019", 3),
    @("Sergipe",  "Número médio de moradores", "01/01/2022", 2.9),
    @("Sergipe",  "Número médio de moradores", "01/01/2023", 2.8),
    @("Sergipe",  "Número médio de moradores", "01/01/2024", 2.8)
)

$startRow = 8
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
}

# Drop the explicit Text number format again now that the values are
# safely stored as text, so the cells don't pick up a style index that
# the original workbook didn't have.
$ws.Range("C8:C22").Style = "Normal"
